# Natmi following Dr Hou advice
# Update NATMI LR-pair metrics for Igf1-Igf1r after the ligand/receptor-expressing-cell
# counts (columns E and K) move from 1 to 3, recomputing every dependent total and
# specificity column (G/H/I/J, M/N/O/P, Q/R/S/T) for rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "E" = 3; "G" = 1.678101666666667; "H" = 5.034305; "I" = 0.00557042851373107; "J" = 0.005570428513731072; "K" = 3; "M" = 14.561928; "N" = 43.685784; "O" = 0.3501777048818433; "P" = 0.3501777048818433; "Q" = 24.43639564668; "R" = 219.92756082012; "S" = 0.001950639872146724; "T" = 0.001950639872146724 }
    3 = @{ "E" = 3; "G" = 1.678101666666667; "H" = 5.034305; "I" = 0.00557042851373107; "J" = 0.005570428513731072; "K" = 3; "M" = 14.40015733333333; "N" = 43.200472; "O" = 0.3462875276490937; "P" = 0.3462875276490937; "Q" = 24.16492802132889; "R" = 217.48435219196; "S" = 0.001928969917965948; "T" = 0.001928969917965949 }
    4 = @{ "E" = 3; "G" = 1.678101666666667; "H" = 5.034305; "I" = 0.00557042851373107; "J" = 0.005570428513731072; "K" = 3; "M" = 2.886742333333333; "N" = 8.660226999999999; "O" = 0.06941888497676431; "P" = 0.06941888497676431; "Q" = 4.844247120803888; "R" = 43.598224087235; "S" = 0.0003866929362659853; "T" = 0.0003866929362659855 }
    5 = @{ "E" = 3; "G" = 1.678101666666667; "H" = 5.034305; "I" = 0.00557042851373107; "J" = 0.005570428513731072; "K" = 3; "M" = 9.735567333333334; "N" = 29.206702; "O" = 0.2341158824922987; "P" = 0.2341158824922987; "Q" = 16.33727176801222; "R" = 147.03544591211; "S" = 0.001304125787352413; "T" = 0.001304125787352414 }
    6 = @{ "E" = 3; "G" = 242.032496; "H" = 726.097488; "I" = 0.8034225480783751; "J" = 0.8034225480783752; "K" = 3; "M" = 14.561928; "N" = 43.685784; "O" = 0.3501777048818433; "P" = 0.3501777048818433; "Q" = 3524.459780412288; "R" = 31720.13802371059; "S" = 0.2813406639364078; "T" = 0.2813406639364078 }
    7 = @{ "E" = 3; "G" = 242.032496; "H" = 726.097488; "I" = 0.8034225480783751; "J" = 0.8034225480783752; "K" = 3; "M" = 14.40015733333333; "N" = 43.200472; "O" = 0.3462875276490937; "P" = 0.3462875276490937; "Q" = 3485.306022179371; "R" = 31367.75419961434; "S" = 0.2782152078315956; "T" = 0.2782152078315957 }
    8 = @{ "E" = 3; "G" = 242.032496; "H" = 726.097488; "I" = 0.8034225480783751; "J" = 0.8034225480783752; "K" = 3; "M" = 2.886742333333333; "N" = 8.660226999999999; "O" = 0.06941888497676431; "P" = 0.06941888497676431; "Q" = 698.6854522455305; "R" = 6288.169070209775; "S" = 0.05577269745279161; "T" = 0.05577269745279162 }
    9 = @{ "E" = 3; "G" = 242.032496; "H" = 726.097488; "I" = 0.8034225480783751; "J" = 0.8034225480783752; "K" = 3; "M" = 9.735567333333334; "N" = 29.206702; "O" = 0.2341158824922987; "P" = 0.2341158824922987; "Q" = 2356.323661662731; "R" = 21206.91295496458; "S" = 0.1880939788575801; "T" = 0.1880939788575801 }
    10 = @{ "E" = 3; "G" = 54.59360333333333; "H" = 163.78081; "I" = 0.1812224912924368; "J" = 0.1812224912924368; "K" = 3; "M" = 14.561928; "N" = 43.685784; "O" = 0.3501777048818433; "P" = 0.3501777048818433; "Q" = 794.98812100056; "R" = 7154.89308900504; "S" = 0.06346007607375534; "T" = 0.06346007607375535 }
    11 = @{ "E" = 3; "G" = 54.59360333333333; "H" = 163.78081; "I" = 0.1812224912924368; "J" = 0.1812224912924368; "K" = 3; "M" = 14.40015733333333; "N" = 43.200472; "O" = 0.3462875276490937; "P" = 0.3462875276490937; "Q" = 786.1564773935912; "R" = 7075.408296542321; "S" = 0.06275508846406734; "T" = 0.06275508846406735 }
    12 = @{ "E" = 3; "G" = 54.59360333333333; "H" = 163.78081; "I" = 0.1812224912924368; "J" = 0.1812224912924368; "K" = 3; "M" = 2.886742333333333; "N" = 8.660226999999999; "O" = 0.06941888497676431; "P" = 0.06941888497676431; "Q" = 157.5976658715411; "R" = 1418.37899284387; "S" = 0.01258026327823234; "T" = 0.01258026327823234 }
    13 = @{ "E" = 3; "G" = 54.59360333333333; "H" = 163.78081; "I" = 0.1812224912924368; "J" = 0.1812224912924368; "K" = 3; "M" = 9.735567333333334; "N" = 29.206702; "O" = 0.2341158824922987; "P" = 0.2341158824922987; "Q" = 531.4997012209578; "R" = 4783.49731098862; "S" = 0.04242706347638175; "T" = 0.04242706347638176 }
    14 = @{ "E" = 3; "G" = 2.947608; "H" = 8.842824; "I" = 0.009784532115456939; "J" = 0.009784532115456941; "K" = 3; "M" = 14.561928; "N" = 43.685784; "O" = 0.3501777048818433; "P" = 0.3501777048818433; "Q" = 42.922855468224; "R" = 386.305699214016; "S" = 0.003426324999533398; "T" = 0.003426324999533398 }
    15 = @{ "E" = 3; "G" = 2.947608; "H" = 8.842824; "I" = 0.009784532115456939; "J" = 0.009784532115456941; "K" = 3; "M" = 14.40015733333333; "N" = 43.200472; "O" = 0.3462875276490937; "P" = 0.3462875276490937; "Q" = 42.44601895699201; "R" = 382.0141706129281; "S" = 0.00338826143546474; "T" = 0.003388261435464741 }
    16 = @{ "E" = 3; "G" = 2.947608; "H" = 8.842824; "I" = 0.009784532115456939; "J" = 0.009784532115456941; "K" = 3; "M" = 2.886742333333333; "N" = 8.660226999999999; "O" = 0.06941888497676431; "P" = 0.06941888497676431; "Q" = 8.508984795671999; "R" = 76.58086316104799; "S" = 0.0006792313094743616; "T" = 0.0006792313094743618 }
    17 = @{ "E" = 3; "G" = 2.947608; "H" = 8.842824; "I" = 0.009784532115456939; "J" = 0.009784532115456941; "K" = 3; "M" = 9.735567333333334; "N" = 29.206702; "O" = 0.2341158824922987; "P" = 0.2341158824922987; "Q" = 28.696636156272; "R" = 258.269725406448; "S" = 0.00229071437098444; "T" = 0.00229071437098444 }
}

foreach ($rowNum in $rowData.Keys) {
    $cols = $rowData[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value2 = $cols[$col]
    }
}

